$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.426.34'
$ws.Range("E2").Value = '  +4.18%  '
$ws.Range("D3").Value = '3.499.48'
$ws.Range("E3").Value = '  +3.88%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'585.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.81%  '
$ws.Range("D6").Value = "'147.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.56%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +1.47%  '
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("E10").Value = '  +4.54%  '
$ws.Range("E11").Value = '  +4.82%  '
$ws.Range("D12").Value = '4.099.70'
$ws.Range("E12").Value = '  +3.95%  '
$ws.Range("D13").Value = "'29.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.12%  '
$ws.Range("D15").Value = '3.501.17'
$ws.Range("E15").Value = '  +3.98%  '
$ws.Range("E16").Value = '  +4.40%  '
$ws.Range("D17").Value = '63.462.21'
$ws.Range("E17").Value = '  +4.10%  '
$ws.Range("E18").Value = '  +3.21%  '
$ws.Range("E19").Value = '  +5.52%  '
$ws.Range("E20").Value = '  +7.11%  '
$ws.Range("D21").Value = "'395.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.64%  '
$ws.Range("E22").Value = '  +3.30%  '
$ws.Range("D23").Value = "'75.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = "'0.0000120"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.25%  '
$ws.Range("D26").Value = '3.643.07'
$ws.Range("E26").Value = '  +3.92%  '
$ws.Range("E27").Value = '  -1.12%  '
$ws.Range("D28").Value = "'7.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.78%  '
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E31").Value = '  +2.95%  '
$ws.Range("D32").Value = "'1.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.82%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  +4.18%  '
$ws.Range("B35").Value = 'EnergySwap'
$ws.Range("C35").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D35").Value = "'32.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +28.41%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").Value = "'7.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.02%  '
$ws.Range("D37").Value = "'5.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.86%  '
$ws.Range("D38").Value = "'172.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.94%  '
$ws.Range("E39").Value = '  +9.43%  '
$ws.Range("D40").Value = '3.537.75'
$ws.Range("E40").Value = '  +3.96%  '
$ws.Range("D41").Value = "'0.0770"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.46%  '
$ws.Range("D42").Value = "'0.803"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.96%  '
$ws.Range("E43").Value = '  +7.94%  '
$ws.Range("D44").Value = "'4.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.32%  '
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("E46").Value = '  +10.31%  '
$ws.Range("D47").Value = '2.612.26'
$ws.Range("E47").Value = '  +6.52%  '
$ws.Range("D48").Value = "'23.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.35%  '
$ws.Range("E49").Value = '  +13.50%  '
$ws.Range("D51").Value = "'0.0270"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.10%  '
